$d = $word.ActiveDocument

$old1 = "1. Origem celular das enzimas: origem das enzimas, diferenciação entre enzimas intra e extracelulares, importância fisiológica e introdução ao mercado mundial de enzimas;2. Estrutura versus propriedades e mecanismos de ação das enzimas: estruturas tridimensionais e sua determinação, importância da estrutura terciária na atividade catalítica, ação catalítica de proteases, glicosidases e oxido-redutases;3. Purificação de enzimas recombinante: métodos de produção, métodos de extração de enzimas, métodos preliminaries de purificação, carreadores enzimáticos, métodos de separação por afinidade;4. Análises massivas utilizadas na prospecção de enzimas de interesse;5. Enzimas imobilizadas: formas de imobilização e aplicações de sistemas imobilizados.6. Métodos utilizados no melhoramento de enzimas de interesse (desenho racional versus evolução direta);7. Enzimas em cosméticos;8. Aplicações de enzimas na indústria: uso de enzimas em detergentes, no processamento do amido, na indústria alimentícia, na indústria têxtil, na síntese de fármacos e na indústria de celulose e papel."
$new1 = "1. Origem celular das enzimas: origem das enzimas, diferenciação entre enzimas intra e extracelulares, importância fisiológica e introdução ao mercado mundial de enzimas;^l2. Estrutura versus propriedades e mecanismos de ação das enzimas: estruturas tridimensionais e sua determinação, importância da estrutura terciária na atividade catalítica, ação catalítica de proteases, glicosidases e oxido-redutases;^l3. Purificação de enzimas recombinante: métodos de produção, métodos de extração de enzimas, métodos preliminaries de purificação, carreadores enzimáticos, métodos de separação por afinidade;^l4. Análises massivas utilizadas na prospecção de enzimas de interesse;^l5. Enzimas imobilizadas: formas de imobilização e aplicações de sistemas imobilizados.^l6. Métodos utilizados no melhoramento de enzimas de interesse (desenho racional versus evolução direta);^l7. Enzimas em cosméticos;^l8. Aplicações de enzimas na indústria: uso de enzimas em detergentes, no processamento do amido, na indústria alimentícia, na indústria têxtil, na síntese de fármacos e na indústria de celulose e papel."
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Output "Replace 1: $found1"

$old2 = "1. Cellular origin of enzymes: origins of enzymes, differentiation between intra- and extracellular enzymes, physiological importance, and introduction to the global enzyme market.2. Structure versus properties and mechanisms of enzyme action: three-dimensional structures and their determination, importance of tertiary structure in catalytic activity, catalytic action of proteases, glycosidases, and oxidoreductases.3. Purification of recombinant enzymes: production methods, enzyme extraction methods, preliminary purification methods, enzyme carriers, affinity separation methods.4. Mass analysis techniques used in enzyme prospecting for interest.5. Immobilized enzymes: immobilization methods and applications of immobilized systems.6. Methods used in enzyme improvement (rational design versus directed evolution).7. Enzymes in cosmetics.8. Applications of enzymes in industry: use of enzymes in detergents, starch processing, food industry, textile industry, drug synthesis, and pulp and paper industry."
$new2 = "1. Cellular origin of enzymes: origins of enzymes, differentiation between intra- and extracellular enzymes, physiological importance, and introduction to the global enzyme market.^l2. Structure versus properties and mechanisms of enzyme action: three-dimensional structures and their determination, importance of tertiary structure in catalytic activity, catalytic action of proteases, glycosidases, and oxidoreductases.^l3. Purification of recombinant enzymes: production methods, enzyme extraction methods, preliminary purification methods, enzyme carriers, affinity separation methods.^l4. Mass analysis techniques used in enzyme prospecting for interest.^l5. Immobilized enzymes: immobilization methods and applications of immobilized systems.^l6. Methods used in enzyme improvement (rational design versus directed evolution).^l7. Enzymes in cosmetics.^l8. Applications of enzymes in industry: use of enzymes in detergents, starch processing, food industry, textile industry, drug synthesis, and pulp and paper industry."
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Output "Replace 2: $found2"

$old3 = "1. Said, S., Pietro, R.C.L. (2002). Enzimas de interesse industrial e biotecnológico. Eventos Editora, Rio de Janeiro.2. Bon, E.S., Ferrara M.A., Corvo M.L. (Eds.) Enzimas em Biotecnologia - Produção, aplicação e mercado, Rio de Janeiro: Editora Interciêcnia, 2008.3.Voet, D., Voet, J., Pratt, C.W. Fundamentos de Bioquímica. Porto Alegre: Editora ARTMED, 2000.4. Walker, J.M., Rapley, R,. Molecular Biomethods Handbook. Humana Press, 2008.5. Copeland, R.A. Enzymes, 3rd Edition. Wiley-Blackwell, 2023."
$new3 = "1. Said, S., Pietro, R.C.L. (2002). Enzimas de interesse industrial e biotecnológico. Eventos Editora, Rio de Janeiro.^l2. Bon, E.S., Ferrara M.A., Corvo M.L. (Eds.) Enzimas em Biotecnologia - Produção, aplicação e mercado, Rio de Janeiro: Editora Interciêcnia, 2008.^l3.Voet, D., Voet, J., Pratt, C.W. Fundamentos de Bioquímica. Porto Alegre: Editora ARTMED, 2000.^l4. Walker, J.M., Rapley, R,. Molecular Biomethods Handbook. Humana Press, 2008.^l5. Copeland, R.A. Enzymes, 3rd Edition. Wiley-Blackwell, 2023."
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Output "Replace 3: $found3"
